# Insert a new weekly price record for "Acelga" (Vega Monumental Concepción)
# as row 238, shifting all subsequent rows down by one (old row 238 becomes
# row 239, ..., old row 321 becomes row 322).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 238 (existing rows 238..321 shift to 239..322)
$ws.Rows.Item(238).Insert()

# Populate the new row 238 with the new record's data
$ws.Cells.Item(238, 1).Value = 11
$ws.Cells.Item(238, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(238, 3).Value = "Bíobío"
$ws.Cells.Item(238, 4).Value = 44876
$ws.Cells.Item(238, 5).Value = 8
$ws.Cells.Item(238, 6).Value = 100112009
$ws.Cells.Item(238, 7).Value = "Acelga"
$ws.Cells.Item(238, 8).Value = "Sin especificar"
$ws.Cells.Item(238, 9).Value = "Primera"
$ws.Cells.Item(238, 10).Value = 250
$ws.Cells.Item(238, 11).Value = 650
$ws.Cells.Item(238, 12).Value = 700
$ws.Cells.Item(238, 13).Value = 670
$ws.Cells.Item(238, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(238, 15).Value = "Región de Ñuble"
$ws.Cells.Item(238, 16).Value = 670
$ws.Cells.Item(238, 17).Value = 1
$ws.Cells.Item(238, 18).Value = "Hortaliza"
